$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the black separator formatting currently on row 4 (A4:E4) down to the
# new separator row 7, extended to span A:G, before row 4's contents are
# overwritten with real data below.
$ws.Range("A4:E4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 4 new test case values written first (matches shared-string order) ---
$ws.Range("C4").Value = "Tim"
$ws.Range("D4").Value = "Bob"

# Header additions
$ws.Range("F1").Value = "nullValue"

# Execute flags on existing rows flipped to "No"
$ws.Range("B2").Value = "No"
$ws.Range("F2").Value = "No"
$ws.Range("B3").Value = "No"
$ws.Range("F3").Value = "No"

# Row 4 test name
$ws.Range("A4").Value = "AddCustomerError1"

# Header for description column
$ws.Range("G1").Value = "Desc"

# Descriptions
$ws.Range("G2").Value = "Adds Customer No Errors"
$ws.Range("G3").Value = "Adds Customer No Errors"

# Remaining Row 4 values
$ws.Range("B4").Value = "Yes"
$ws.Range("F4").Value = "PostCode"
$ws.Range("G4").Value = "Mandatory Field Error PostCode pops up"

# Row 5 - New test case: AddCustomerError2
$ws.Range("C5").Value = "John"
$ws.Range("E5").Value = 1011
$ws.Range("F5").Value = "LastName"
$ws.Range("G5").Value = "Mandatory Field Error LastName pops up"
$ws.Range("A5").Value = "AddCustomerError2"
$ws.Range("B5").Value = "Yes"

# Row 6 - New test case: AddCustomerError3
$ws.Range("A6").Value = "AddCustomerError3"
$ws.Range("B6").Value = "Yes"
$ws.Range("D6").Value = "Gary"
$ws.Range("E6").Value = 1011
$ws.Range("G6").Value = "Mandatory Field Error FirstName pops up"

# Clear the old black-fill styling from row 4 (it now holds data, no style)
$ws.Range("A4:D4").Style = "Normal"
$ws.Range("F4:G4").Style = "Normal"
# E4 has no data in the new layout - remove it entirely (was part of the old
# black separator row, now unused since the separator moved to row 7)
$ws.Range("E4").Clear()

# Column widths, best-fit for the new longer values in columns A and G
# (target stored widths are 17.7109375 / 37.28515625 characters; the inputs
# below are chosen so that, after the host's internal width quantization,
# the saved width lands as close as possible to those values).
$ws.Columns.Item(1).ColumnWidth = 16.833333333333336
$ws.Columns.Item(7).ColumnWidth = 36.5

$ws.Range("I17").Select()
